# Update "10 min frame" forecast report to a new frame (values in column B
# and the frame length in column C, rows 2-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = 806;   C = 13 },
    @{ Row = 3;  B = 338;   C = 13 },
    @{ Row = 4;  B = 7;     C = 13 },
    @{ Row = 5;  B = 144;   C = 13 },
    @{ Row = 6;  B = 48;    C = 13 },
    @{ Row = 7;  B = 110;   C = 13 },
    @{ Row = 8;  B = 29;    C = 13 },
    @{ Row = 9;  B = 71;    C = 13 },
    @{ Row = 10; B = 614;   C = 13 },
    @{ Row = 11; B = 309;   C = 13 },
    @{ Row = 12; B = 789.1; C = 13 },
    @{ Row = 13; B = 483;   C = 13 },
    @{ Row = 14; B = 719;   C = 13 },
    @{ Row = 15; B = 225;   C = 13 },
    @{ Row = 16; B = 151;   C = 13 },
    @{ Row = 17; B = 107;   C = 13 },
    @{ Row = 18; B = 5;     C = 13 },
    @{ Row = 19; B = 13;    C = 13 },
    @{ Row = 20; B = 61;    C = 13 }
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
    $ws.Cells.Item($entry.Row, 3).Value = $entry.C
}
